$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new columns before column F (old F..Q shift to I..T)
$ws.Columns("F:H").Insert()

# Rename / set header row text
$ws.Range("E1").Value = "training_type"
$ws.Range("F1").Value = "pgd_train_eps"
$ws.Range("G1").Value = "pgd_train_eta"
$ws.Range("H1").Value = "pgd_train_num_iter"

# Add three new data rows (25-27) for PGD hyperparameter-search runs
$ws.Cells.Item(25, 1).Value = 23
$ws.Cells.Item(25, 2).Value = 2
$ws.Cells.Item(25, 3).Value = 0
$ws.Cells.Item(25, 4).Value = 0.003
$ws.Cells.Item(25, 5).Value = "PGD"
$ws.Cells.Item(25, 6).Value = 0.1
$ws.Cells.Item(25, 7).Value = 0.1
$ws.Cells.Item(25, 8).Value = 3
$ws.Cells.Item(25, 9).Value = "<function relu at 0x118b969d8>"
$ws.Cells.Item(25, 10).Value = 0.8787999749183655
$ws.Cells.Item(25, 11).Value = 0.01040000002831221
$ws.Cells.Item(25, 12).Value = 0.001300000003539026
$ws.Cells.Item(25, 13).Value = 0.4310351014137268
$ws.Cells.Item(25, 14).Value = 8.296195983886719
$ws.Cells.Item(25, 15).Value = 0.01040000002831221
$ws.Cells.Item(25, 16).Value = "logs/results_278.log"
$ws.Cells.Item(25, 17).Value = "weights/model_278.ckpt"
$ws.Cells.Item(25, 18).Value = "tb/278/robust"

$ws.Cells.Item(26, 1).Value = 24
$ws.Cells.Item(26, 2).Value = 2
$ws.Cells.Item(26, 3).Value = 0
$ws.Cells.Item(26, 4).Value = 0.003
$ws.Cells.Item(26, 5).Value = "PGD"
$ws.Cells.Item(26, 6).Value = 0.1
$ws.Cells.Item(26, 7).Value = 0.1
$ws.Cells.Item(26, 8).Value = 3
$ws.Cells.Item(26, 9).Value = "<function relu at 0x121b0f9d8>"
$ws.Cells.Item(26, 10).Value = 0.8522999882698059
$ws.Cells.Item(26, 11).Value = 0.002899999963119626
$ws.Cells.Item(26, 12).Value = 0
$ws.Cells.Item(26, 13).Value = 0.5192863941192627
$ws.Cells.Item(26, 14).Value = 9.538139343261719
$ws.Cells.Item(26, 15).Value = 0.002899999963119626
$ws.Cells.Item(26, 16).Value = "logs/results_279.log"
$ws.Cells.Item(26, 17).Value = "weights/model_279.ckpt"
$ws.Cells.Item(26, 18).Value = "tb/279/robust"

$ws.Cells.Item(27, 1).Value = 25
$ws.Cells.Item(27, 2).Value = 2
$ws.Cells.Item(27, 3).Value = 0
$ws.Cells.Item(27, 4).Value = 0.003
$ws.Cells.Item(27, 5).Value = "PGD"
$ws.Cells.Item(27, 6).Value = 0.1
$ws.Cells.Item(27, 7).Value = 0.1
$ws.Cells.Item(27, 8).Value = 3
$ws.Cells.Item(27, 9).Value = "<function relu at 0x121b0f9d8>"
$ws.Cells.Item(27, 10).Value = 0.8676999807357788
$ws.Cells.Item(27, 11).Value = 0.004999999888241291
$ws.Cells.Item(27, 12).Value = 0.00009999999747378752
$ws.Cells.Item(27, 13).Value = 0.4790646433830261
$ws.Cells.Item(27, 14).Value = 9.238405227661133
$ws.Cells.Item(27, 15).Value = 0.004999999888241291
$ws.Cells.Item(27, 16).Value = "logs/results_279.log"
$ws.Cells.Item(27, 17).Value = "weights/model_279.ckpt"
$ws.Cells.Item(27, 18).Value = "tb/279/robust"
